$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix K4: make it an explicit (non-shared) formula ---
$ws.Range("K4").Formula = "=J4*0.15"

# --- K10: updated value ---
$ws.Range("K10").Value = 9.8428938354397566

# --- New row 11 addition: L11 formula (depends on I18, filled in later) ---
$ws.Range("L11").Formula = "=J11*I18"

# --- New row 12: Cephalosporins data ---
$ws.Range("A12").Value = "Cephalosporins"
$ws.Range("B12").Value = 540
$ws.Range("C12").Value = 1000
$ws.Range("D12").Value = 433.33330000000001
$ws.Range("E12").Value = 500
$ws.Range("J12").Formula = "=AVERAGE(B12:E12)"
$ws.Range("K12").Formula = "=J12*0.15"
$ws.Range("L12").Value = 18964.994444444445
$ws.Range("M12").Formula = "=K12/L12"

# --- Row 16: widen average range to include new row ---
$ws.Range("I16").Formula = "=AVERAGE(M2:M12)"

# --- New row 17: average_price label + formula ---
$ws.Range("L17").Value = "average_price"
$ws.Range("M17").Formula = "=AVERAGE(K2:K12)"

# --- Row 18 formula stays the same text, but recalculates with new I16/L10 ---
$ws.Range("I18").Formula = "=L10*I16"

# --- Sheet view / selection update ---
$ws.Range("K9").Select()
